# Apply cryptocurrency price/volume updates as described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "62.273.55"
$cell.Style = "Normal"
$ws.Range("E2").Value = "  -0.64%  "

# Row 3
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "2.457.11"
$cell.Style = "Normal"
$ws.Range("E3").Value = "  +0.62%  "

# Row 4
$ws.Range("E4").Value = "  +0.00%  "

# Row 5
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "574.58"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  +0.96%  "

# Row 6
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "144.05"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  -0.73%  "

# Row 7
$ws.Range("E7").Value = "  +0.08%  "

# Row 8
$ws.Range("E8").Value = "  +0.09%  "

# Row 9
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "2.453.50"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  +0.67%  "

# Row 10
$ws.Range("E10").Value = "  -0.82%  "

# Row 11
$ws.Range("E11").Value = "  +2.13%  "

# Row 12
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "5.22"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  -0.10%  "

# Row 13
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "0.345"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  -2.97%  "

# Row 14
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "26.34"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  -2.55%  "

# Row 15
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "0.0000174"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  -1.09%  "

# Row 16
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "2.894.98"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  +0.36%  "

# Row 17
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "62.181.95"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  -0.58%  "

# Row 18
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "2.450.24"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  +0.56%  "

# Row 19
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "10.86"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  -3.39%  "

# Row 20
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "7.16"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  -1.01%  "

# Row 21
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "328.28"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  +0.23%  "

# Row 22
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "4.13"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  -0.98%  "

# Row 23
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "1.95"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  -6.11%  "

# Row 24
$ws.Range("E24").Value = "  +0.02%  "

# Row 25
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "65.77"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  +0.87%  "

# Row 26
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "9.21"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  +1.75%  "

# Row 27
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "601.15"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  -3.15%  "

# Row 28
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "2.581.85"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  +1.00%  "

# Row 29
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "0.0₃0961"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  -3.62%  "

# Row 30
$ws.Range("E30").Value = "  -0.09%  "

# Row 31
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "1.43"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  -4.62%  "

# Row 32
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "8.00"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  -1.81%  "

# Row 33
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "1.87"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  -0.29%  "

# Row 34
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "0.138"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  -0.36%  "

# Row 35
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "4.91"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  -4.38%  "

# Row 36
$ws.Range("E36").Value = "  +0.33%  "

# Row 37
$ws.Range("E37").Value = "  -3.49%  "

# Row 38
$ws.Range("E38").Value = "  +0.47%  "

# Row 39
$ws.Range("B39").Value = "Monero"
$ws.Range("C39").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "150.86"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  +2.85%  "

# Row 40
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "5.37"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  +0.68%  "

# Row 41
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "18.43"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  -1.84%  "

# Row 42
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "1.73"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  -2.76%  "

# Row 43
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "42.66"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  +1.48%  "

# Row 44
$ws.Range("E44").Value = "  +0.00%  "

# Row 45
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "2.51"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  -2.82%  "

# Row 46
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "142.77"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  -2.26%  "

# Row 47
$ws.Range("E47").Value = "  -3.25%  "

# Row 48
$ws.Range("E48").Value = "  +18.04%  "

# Row 49
$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "0.606"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  +1.65%  "

# Row 50
$ws.Range("B50").Value = "Hedera"
$ws.Range("C50").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "0.0525"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  -0.54%  "

# Row 51
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "19.69"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  -4.79%  "
